$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# --- Update time_taken (column F) timestamps on the "data" sheet ---
$timestamps = @{
    2 = "2021-10-05 14:33:55.837768"
    3 = "2021-10-05 14:33:55.837776"
    4 = "2021-10-05 14:33:55.837779"
    5 = "2021-10-05 14:33:55.837781"
    6 = "2021-10-05 14:33:55.837784"
    7 = "2021-10-05 14:33:55.837787"
    8 = "2021-10-05 14:33:55.837790"
    9 = "2021-10-05 14:33:55.837792"
    10 = "2021-10-05 14:33:55.837795"
    11 = "2021-10-05 14:33:55.837797"
    12 = "2021-10-05 14:33:55.837800"
    13 = "2021-10-05 14:33:55.837802"
    14 = "2021-10-05 14:33:55.837805"
    15 = "2021-10-05 14:33:55.837808"
    16 = "2021-10-05 14:33:55.837810"
    17 = "2021-10-05 14:33:55.837812"
    18 = "2021-10-05 14:33:55.837815"
    19 = "2021-10-05 14:33:55.837818"
    20 = "2021-10-05 14:33:55.837820"
    21 = "2021-10-05 14:33:55.837823"
    22 = "2021-10-05 14:33:55.837826"
    23 = "2021-10-05 14:33:55.837828"
    24 = "2021-10-05 14:33:55.837831"
    25 = "2021-10-05 14:33:55.837833"
    26 = "2021-10-05 14:33:55.837836"
    27 = "2021-10-05 14:33:55.837839"
    28 = "2021-10-05 14:33:55.837841"
    29 = "2021-10-05 14:33:55.837844"
    30 = "2021-10-05 14:33:55.837846"
    31 = "2021-10-05 14:33:55.837848"
    32 = "2021-10-05 14:33:55.837851"
    33 = "2021-10-05 14:33:55.837853"
    34 = "2021-10-05 14:33:55.837856"
    35 = "2021-10-05 14:33:55.837859"
    36 = "2021-10-05 14:33:55.837861"
    37 = "2021-10-05 14:33:55.837864"
    38 = "2021-10-05 14:33:55.837866"
    39 = "2021-10-05 14:33:55.837869"
    40 = "2021-10-05 14:33:55.837871"
    41 = "2021-10-05 14:33:55.837873"
    42 = "2021-10-05 14:33:55.837877"
    43 = "2021-10-05 14:33:55.837879"
    44 = "2021-10-05 14:33:55.837882"
    45 = "2021-10-05 14:33:55.837884"
    46 = "2021-10-05 14:33:55.837887"
    47 = "2021-10-05 14:33:55.837889"
    48 = "2021-10-05 14:33:55.837891"
    49 = "2021-10-05 14:33:55.837894"
    50 = "2021-10-05 14:33:55.837896"
    51 = "2021-10-05 14:33:55.837899"
    52 = "2021-10-05 14:33:55.837901"
}
foreach ($row in $timestamps.Keys) {
    $dataWs.Cells.Item([int]$row, 6).Value = $timestamps[$row]
}

# --- Add the new "metadata" worksheet right after "data" ---
$metaWs = $wb.Worksheets.Add($null, $dataWs)
$metaWs.Name = "metadata"

# Header row (B1:G1) -- styled like the "data" sheet header (bold + border via copy/paste of formats)
$metaWs.Range("B1").Value = "data_name"
$metaWs.Range("C1").Value = "data_id"
$metaWs.Range("D1").Value = "data_version"
$metaWs.Range("E1").Value = "data_version_created"
$metaWs.Range("F1").Value = "panel_query_time"
$metaWs.Range("G1").Value = "panel_get_request"

# Data row 2
$metaWs.Range("A2").Value = 0
$metaWs.Range("B2").Value = "Hair disorders"
$metaWs.Range("C2").Value = 3269
$metaWs.Range("D2").Value = "0.46"
$metaWs.Range("E2").Value = "2021-07-03T07:21:44.767485Z"
$metaWs.Range("F2").Value = "2021-10-05 14:33:55.834016"
$metaWs.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3269/?format=json"

# Copy the existing bold/bordered header style (cellXfs index 1, used by the "data" header row)
# onto the new header cells and the A2 index cell, reusing the style instead of creating a new one.
$dataWs.Range("B1").Copy()
$metaWs.Range("B1:G1").PasteSpecial(-4122)
$dataWs.Range("A2").Copy()
$metaWs.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
